# [BOM] Update bill of materials - Replace feedback resistors
#
# Row 7 of the "BLM_2017" BOM sheet holds the R2/R3 feedback resistors
# (Designator "R2, R3"). The new parts get different manufacturer /
# supplier part numbers and a higher unit price, which ripples into the
# line subtotal and the J9 grand-total formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BLM_2017")

# Update the part numbers. Write H7 before D7 so the new shared strings
# land in the same table order as the source edit.
$ws.Range("H7").Value = "P3.3AJCT-ND"
$ws.Range("D7").Value = "ERJ-3RQF3R3V"

# Writing new text into these cells drops their original quote-prefixed
# text format; restore it from a same-styled neighbour cell.
$ws.Range("H6").Copy() | Out-Null
$ws.Range("H7").PasteSpecial(-4122) | Out-Null
$ws.Range("D6").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null

# New unit price / subtotal for the R2, R3 line.
$ws.Range("I7").Value = 0.42
$ws.Range("J7").Value = 0.84

$wb.Save()
